# The "Requisitos" bullet paragraph lists three course-code lines, each its
# own run ending in a manual line break (<w:br/>). The edit reorders them by
# moving the "LOM3246 - Técnicas de Caracterização de Materiais ..." line
# from the top of the list to the bottom (after the "LOM3016 ..." line),
# while leaving "LOB1021 ..." and "LOM3016 ..." untouched.

$d = $word.ActiveDocument

$lom3246Line = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)" + [char]11
$lom3016Line = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)" + [char]11

# 1) Insert a fresh copy of the LOM3246 line (text + break) right after the
#    LOM3016 line's break. InsertAfter creates it as its own run, matching
#    the target OOXML (three sibling <w:r> elements).
$lom3016Range = $d.Content
$found = $lom3016Range.Find.Execute($lom3016Line, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $lom3016Range.Collapse(0)
    $lom3016Range.InsertAfter($lom3246Line)
}

# 2) Remove the original LOM3246 line (text + break) from the top of the list.
$lom3246Range = $d.Content
$found2 = $lom3246Range.Find.Execute($lom3246Line, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $lom3246Range.Delete()
}
